# Updates the cryptos list (Price column D, Volume(1h) column E) to match
# the latest scraped values, per commit 'Updated cryptos list ... with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are plain text in the sheet (no numeric format applied).
# Assigning a numeric-looking string via .Value would silently convert the cell to a
# true Number (losing trailing zeros / introducing float rounding), so we force the
# cell through a quoted-text Formula, flip its stored data type back to Text, then
# restore the 'Normal' cell style so no stray number-format/quote-prefix style sticks.
function Set-TextValue($cell, $text) {
    $cell.Formula = "'" + $text
    $cell.DataTypeToText()
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "57.881.68"
$ws.Cells.Item(2, 5).Value = "  +2.88%  "
Set-TextValue $ws.Cells.Item(3, 4) "3.068.73"
$ws.Cells.Item(3, 5).Value = "  +2.39%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
Set-TextValue $ws.Cells.Item(5, 4) "516.92"
$ws.Cells.Item(5, 5).Value = "  +1.90%  "
Set-TextValue $ws.Cells.Item(6, 4) "142.37"
$ws.Cells.Item(6, 5).Value = "  +2.79%  "
$ws.Cells.Item(7, 5).Value = "  +0.01%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.436"
$ws.Cells.Item(8, 5).Value = "  +1.90%  "
$ws.Cells.Item(9, 5).Value = "  +2.67%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.107"
$ws.Cells.Item(10, 5).Value = "  +0.48%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.377"
$ws.Cells.Item(11, 5).Value = "  +2.78%  "
Set-TextValue $ws.Cells.Item(12, 4) "3.596.60"
$ws.Cells.Item(12, 5).Value = "  +2.32%  "
$ws.Cells.Item(13, 5).Value = "  +3.17%  "
Set-TextValue $ws.Cells.Item(14, 4) "26.17"
$ws.Cells.Item(14, 5).Value = "  +3.34%  "
$ws.Cells.Item(15, 5).Value = "  +1.30%  "
Set-TextValue $ws.Cells.Item(16, 4) "57.918.73"
$ws.Cells.Item(16, 5).Value = "  +3.04%  "
Set-TextValue $ws.Cells.Item(17, 4) "3.068.32"
$ws.Cells.Item(17, 5).Value = "  +2.29%  "
Set-TextValue $ws.Cells.Item(18, 4) "6.09"
$ws.Cells.Item(18, 5).Value = "  +2.95%  "
Set-TextValue $ws.Cells.Item(19, 4) "12.86"
Set-TextValue $ws.Cells.Item(20, 4) "8.11"
$ws.Cells.Item(20, 5).Value = "  +1.71%  "
Set-TextValue $ws.Cells.Item(21, 4) "333.02"
$ws.Cells.Item(21, 5).Value = "  +0.62%  "
Set-TextValue $ws.Cells.Item(22, 4) "0.999"
$ws.Cells.Item(22, 5).Value = "  -0.02%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.500"
$ws.Cells.Item(23, 5).Value = "  +0.91%  "
Set-TextValue $ws.Cells.Item(24, 4) "65.57"
$ws.Cells.Item(24, 5).Value = "  +1.48%  "
$ws.Cells.Item(25, 5).Value = "  +3.27%  "
Set-TextValue $ws.Cells.Item(26, 4) "1.00"
$ws.Cells.Item(26, 5).Value = "  +0.01%  "
Set-TextValue $ws.Cells.Item(27, 4) "0.0₃0903"
$ws.Cells.Item(27, 5).Value = "  -2.78%  "
Set-TextValue $ws.Cells.Item(28, 4) "6.48"
$ws.Cells.Item(28, 5).Value = "  +2.62%  "
Set-TextValue $ws.Cells.Item(29, 4) "7.27"
$ws.Cells.Item(29, 5).Value = "  +6.22%  "
Set-TextValue $ws.Cells.Item(30, 4) "1.82"
$ws.Cells.Item(30, 5).Value = "  +2.22%  "
$ws.Cells.Item(31, 5).Value = "  +3.20%  "
Set-TextValue $ws.Cells.Item(32, 4) "20.71"
$ws.Cells.Item(32, 5).Value = "  +2.19%  "
Set-TextValue $ws.Cells.Item(33, 4) "154.85"
$ws.Cells.Item(33, 5).Value = "  +1.63%  "
Set-TextValue $ws.Cells.Item(34, 4) "4.55"
$ws.Cells.Item(34, 5).Value = "  +3.15%  "
Set-TextValue $ws.Cells.Item(35, 4) "6.02"
$ws.Cells.Item(35, 5).Value = "  +3.96%  "
Set-TextValue $ws.Cells.Item(36, 4) "26.92"
$ws.Cells.Item(36, 5).Value = "  +1.63%  "
Set-TextValue $ws.Cells.Item(37, 4) "1.27"
$ws.Cells.Item(37, 5).Value = "  +4.11%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.0676"
$ws.Cells.Item(38, 5).Value = "  +2.68%  "
Set-TextValue $ws.Cells.Item(39, 4) "3.111.64"
Set-TextValue $ws.Cells.Item(40, 4) "3.92"
Set-TextValue $ws.Cells.Item(41, 4) "36.44"
$ws.Cells.Item(41, 5).Value = "  +0.18%  "
$ws.Cells.Item(42, 5).Value = "  +0.01%  "
Set-TextValue $ws.Cells.Item(43, 4) "0.656"
$ws.Cells.Item(43, 5).Value = "  +0.35%  "
Set-TextValue $ws.Cells.Item(44, 4) "2.267.18"
$ws.Cells.Item(44, 5).Value = "  +3.61%  "
$ws.Cells.Item(45, 5).Value = "  +8.86%  "
Set-TextValue $ws.Cells.Item(46, 4) "20.84"
$ws.Cells.Item(46, 5).Value = "  +7.28%  "
Set-TextValue $ws.Cells.Item(47, 4) "1.37"
$ws.Cells.Item(47, 5).Value = "  +2.61%  "
$ws.Cells.Item(48, 5).Value = "  +3.21%  "
Set-TextValue $ws.Cells.Item(49, 4) "5.92"
$ws.Cells.Item(49, 5).Value = "  +2.00%  "
Set-TextValue $ws.Cells.Item(50, 4) "0.743"
$ws.Cells.Item(50, 5).Value = "  +10.76%  "
Set-TextValue $ws.Cells.Item(51, 4) "256.94"
$ws.Cells.Item(51, 5).Value = "  +13.14%  "
